$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the "Periodo Mora" values for the first and third data rows
# (1712 <-> 1710), the middle row (1711) stays the same.
$ws.Range("E16").Value = "1710"
$ws.Range("E18").Value = "1712"

# Update "Salario Basico" values for the three data rows
$ws.Range("G16").Value = 781242
$ws.Range("G17").Value = 781242
$ws.Range("G18").Value = 781242

# Adjust column widths (auto-fit widths shrank after the data refresh)
$ws.Columns.Item(2).ColumnWidth = 16.0
$ws.Columns.Item(3).ColumnWidth = 8.0
$ws.Columns.Item(4).ColumnWidth = 28.0
$ws.Columns.Item(5).ColumnWidth = 11.833333333333334
$ws.Columns.Item(6).ColumnWidth = 8.666666666666666
$ws.Columns.Item(7).ColumnWidth = 12.666666666666666
$ws.Columns.Item(8).ColumnWidth = 17.0
$ws.Columns.Item(9).ColumnWidth = 16.0
$ws.Columns.Item(10).ColumnWidth = 13.333333333333334

# The logo shifts left along with the narrower column B - keep its size
# fixed and nudge its position to match.
$logo = $ws.Shapes.Item(1)
$logo.Width = 76.81889763779527
$logo.Height = 48.188976377952756
$logo.Left = 53.59055118110236
